{"js": "// Replace the 25 division-problem answers in the single body table with\n// their new values, in document order (row-major, skipping the blank\n// spacer rows). Positional (index-based) replacement is required because\n// several original answers are duplicated (e.g. \"81\u00f79=9, 0\" appears\n// twice but maps to two different new values), so a text-based find &\n// replace-all would be ambiguous.\nconst replacements = [\n  \"57\u00f79=6, 3\",\n  \"70\u00f72=35, 0\",\n  \"55\u00f72=27, 1\",\n  \"10\u00f72=5, 0\",\n  \"30\u00f76=5, 0\",\n  \"83\u00f75=16, 3\",\n  \"35\u00f75=7, 0\",\n  \"91\u00f77=13, 0\",\n  \"53\u00f73=17, 2\",\n  \"51\u00f72=25, 1\",\n  \"95\u00f79=10, 5\",\n  \"29\u00f74=7, 1\",\n  \"31\u00f74=7, 3\",\n  \"96\u00f75=19, 1\",\n  \"34\u00f76=5, 4\",\n  \"64\u00f72=32, 0\",\n  \"68\u00f72=34, 0\",\n  \"10\u00f75=2, 0\",\n  \"96\u00f73=32, 0\",\n  \"70\u00f75=14, 0\",\n  \"92\u00f79=10, 2\",\n  \"40\u00f76=6, 4\",\n  \"41\u00f77=5, 6\",\n  \"60\u00f72=30, 0\",\n  \"33\u00f77=4, 5\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Gather every cell proxy across all rows in row-major order, then load\n// their current text to know which rows actually hold answers (the\n// spacer rows' paragraphs are empty).\nconst allCells = [];\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  allCells.push(cells);\n}\nawait context.sync();\n\nconst cellList = [];\nfor (const cells of allCells) {\n  for (const cell of cells.items) {\n    cell.load(\"value\");\n    cellList.push(cell);\n  }\n}\nawait context.sync();\n\nlet i = 0;\nfor (const cell of cellList) {\n  const text = cell.value;\n  if (text && text.trim().length > 0) {\n    cell.value = replacements[i];\n    i += 1;\n  }\n}\n\nif (i !== replacements.length) {\n  throw new Error(\n    `Expected to replace ${replacements.length} cells, replaced ${i}`\n  );\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem answers in the single body table with\n# their new values, in document order (row-major, skipping the blank\n# spacer rows). Positional (row/column) replacement is required because\n# several original answers are duplicated (e.g. \"81\u00f79=9, 0\" appears\n# twice but maps to two different new values), so a Find/Replace-All by\n# text would be ambiguous.\n$replacements = @(\n  \"57\u00f79=6, 3\",\n  \"70\u00f72=35, 0\",\n  \"55\u00f72=27, 1\",\n  \"10\u00f72=5, 0\",\n  \"30\u00f76=5, 0\",\n  \"83\u00f75=16, 3\",\n  \"35\u00f75=7, 0\",\n  \"91\u00f77=13, 0\",\n  \"53\u00f73=17, 2\",\n  \"51\u00f72=25, 1\",\n  \"95\u00f79=10, 5\",\n  \"29\u00f74=7, 1\",\n  \"31\u00f74=7, 3\",\n  \"96\u00f75=19, 1\",\n  \"34\u00f76=5, 4\",\n  \"64\u00f72=32, 0\",\n  \"68\u00f72=34, 0\",\n  \"10\u00f75=2, 0\",\n  \"96\u00f73=32, 0\",\n  \"70\u00f75=14, 0\",\n  \"92\u00f79=10, 2\",\n  \"40\u00f76=6, 4\",\n  \"41\u00f77=5, 6\",\n  \"60\u00f72=30, 0\",\n  \"33\u00f77=4, 5\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$i = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $clean = $cell.Range.Text -replace \"[\\r\\a]\", \"\"\n    if ($clean.Length -gt 0) {\n      $cell.Range.Text = $replacements[$i]\n      $i = $i + 1\n    }\n  }\n}\n\nif ($i -ne $replacements.Length) {\n  throw \"Expected to replace $($replacements.Length) cells, replaced $i\"\n}\n"}
